$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A date cells to remain plain text so Excel does not
# auto-convert the dash-separated strings into date serial numbers.
$dateCells = @("A3","A4","A5","A6","A7","A8","A9","A10","A11","A12",
               "A13","A14","A15","A16","A17","A18","A19","A20","A21")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update date strings in column A: replace "/" separators with "-"
$ws.Range("A3").Value  = "28-07-2022"
$ws.Range("A4").Value  = "01-08-2022"
$ws.Range("A5").Value  = "04-08-2022"
$ws.Range("A6").Value  = "08-08-2022"
$ws.Range("A7").Value  = "11-08-2022"
$ws.Range("A8").Value  = "15-08-2022"
$ws.Range("A9").Value  = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# Update attendance counters that changed alongside the date fixes
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("G6").Value = 1

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0
